$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Round row 5 values to 2 decimal places ("custom accuracy")
$row5 = @{
    "B5"  = 13.15
    "C5"  = 9.93
    "D5"  = 0.33
    "E5"  = 27.79
    "F5"  = 23.19
    "G5"  = 9.86
    "H5"  = 42.42
    "I5"  = 15.5
    "J5"  = 7.39
    "K5"  = 10.51
    "L5"  = 11.34
    "M5"  = 12.19
    "N5"  = 3.57
    "O5"  = 9.97
    "P5"  = 14.52
    "Q5"  = 8.13
    "R5"  = 0.38
    "S5"  = 0.3
    "T5"  = 147.02
    "U5"  = 28.17
    "V5"  = 9.57
    "W5"  = 19.22
    "X5"  = 9.86
    "Y5"  = 1.33
    "Z5"  = 20.42
    "AA5" = 8.22
    "AB5" = 7.33
    "AC5" = 9.1
    "AD5" = 11.93
    "AE5" = 0.57
    "AF5" = 38.7
    "AG5" = 5.21
    "AH5" = 11.55
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# 2) Remove the last data row (row 6) -- "데이터 1000개" trims the dataset
$ws.Rows.Item(6).Delete()

# 3) Narrow a handful of columns from width 8 to width 7
#    (ColumnWidth 6.17 chars maps to the stored OOXML width of 7,
#     matching the width already used by the other "7"-wide columns)
$narrowCols = @(6, 11, 21, 23, 32)
foreach ($c in $narrowCols) {
    $ws.Columns.Item($c).ColumnWidth = 6.17
}
